$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date (column C) from 45204 (2023-10-05) to
# 45207 (2023-10-08) for rows 2 through 22.
$ws.Range("C2:C22").Value = 45207
